$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44260
$ws.Cells.Item(2, 10).Value = 33
$ws.Cells.Item(2, 11).Value = 22000
$ws.Cells.Item(2, 12).Value = 23000
$ws.Cells.Item(2, 13).Value = 22545
$ws.Cells.Item(2, 14).Value = '$/saco 30 kilos'
$ws.Cells.Item(2, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(2, 16).Value = 22545
$ws.Cells.Item(2, 17).Value = 1
# Row 4
$ws.Cells.Item(4, 4).Value = 44174
$ws.Cells.Item(4, 11).Value = 19000
$ws.Cells.Item(4, 12).Value = 20000
$ws.Cells.Item(4, 13).Value = 19500
$ws.Cells.Item(4, 15).Value = 'Región de Ñuble'
$ws.Cells.Item(4, 16).Value = 780
# Row 5
$ws.Cells.Item(5, 4).Value = 44258
$ws.Cells.Item(5, 10).Value = 32
$ws.Cells.Item(5, 11).Value = 22000
$ws.Cells.Item(5, 12).Value = 23000
$ws.Cells.Item(5, 13).Value = 22562
$ws.Cells.Item(5, 14).Value = '$/saco 30 kilos'
$ws.Cells.Item(5, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(5, 16).Value = 22562
$ws.Cells.Item(5, 17).Value = 1
# Row 6
$ws.Cells.Item(6, 4).Value = 44498
$ws.Cells.Item(6, 8).Value = 'Sin especificar'
$ws.Cells.Item(6, 10).Value = 120
$ws.Cells.Item(6, 11).Value = 17000
$ws.Cells.Item(6, 12).Value = 18000
$ws.Cells.Item(6, 13).Value = 17500
$ws.Cells.Item(6, 16).Value = 700
# Row 7
$ws.Cells.Item(7, 4).Value = 44223
$ws.Cells.Item(7, 10).Value = 42
$ws.Cells.Item(7, 11).Value = 26000
$ws.Cells.Item(7, 12).Value = 28000
$ws.Cells.Item(7, 13).Value = 27048
$ws.Cells.Item(7, 16).Value = 27048
# Row 8
$ws.Cells.Item(8, 4).Value = 44165
$ws.Cells.Item(8, 10).Value = 42
$ws.Cells.Item(8, 11).Value = 18000
$ws.Cells.Item(8, 12).Value = 19000
$ws.Cells.Item(8, 13).Value = 18595
$ws.Cells.Item(8, 16).Value = 744
# Row 9
$ws.Cells.Item(9, 4).Value = 44167
$ws.Cells.Item(9, 10).Value = 40
$ws.Cells.Item(9, 11).Value = 18000
$ws.Cells.Item(9, 12).Value = 19000
$ws.Cells.Item(9, 13).Value = 18500
$ws.Cells.Item(9, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(9, 15).Value = 'Región del Maule'
$ws.Cells.Item(9, 16).Value = 740
$ws.Cells.Item(9, 17).Value = 25
# Row 10
$ws.Cells.Item(10, 4).Value = 44176
$ws.Cells.Item(10, 11).Value = 20000
$ws.Cells.Item(10, 12).Value = 21000
$ws.Cells.Item(10, 13).Value = 20500
$ws.Cells.Item(10, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(10, 16).Value = 820
# Row 11
$ws.Cells.Item(11, 4).Value = 44161
$ws.Cells.Item(11, 10).Value = 33
$ws.Cells.Item(11, 11).Value = 19000
$ws.Cells.Item(11, 12).Value = 19500
$ws.Cells.Item(11, 13).Value = 19303
$ws.Cells.Item(11, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(11, 16).Value = 772
# Row 12
$ws.Cells.Item(12, 4).Value = 44161
$ws.Cells.Item(12, 10).Value = 34
$ws.Cells.Item(12, 11).Value = 19500
$ws.Cells.Item(12, 13).Value = 19735
$ws.Cells.Item(12, 15).Value = 'Región del Maule'
$ws.Cells.Item(12, 16).Value = 789
# Row 13
$ws.Cells.Item(13, 4).Value = 44162
$ws.Cells.Item(13, 10).Value = 50
$ws.Cells.Item(13, 11).Value = 18500
$ws.Cells.Item(13, 12).Value = 19000
$ws.Cells.Item(13, 13).Value = 18820
$ws.Cells.Item(13, 15).Value = 'Región del Maule'
$ws.Cells.Item(13, 16).Value = 753
# Row 14
$ws.Cells.Item(14, 4).Value = 44159
$ws.Cells.Item(14, 10).Value = 28
$ws.Cells.Item(14, 11).Value = 19000
$ws.Cells.Item(14, 12).Value = 19500
$ws.Cells.Item(14, 13).Value = 19268
$ws.Cells.Item(14, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(14, 16).Value = 771
# Row 15
$ws.Cells.Item(15, 4).Value = 44159
$ws.Cells.Item(15, 10).Value = 56
$ws.Cells.Item(15, 11).Value = 19000
$ws.Cells.Item(15, 12).Value = 20000
$ws.Cells.Item(15, 13).Value = 19464
$ws.Cells.Item(15, 16).Value = 779
# Row 16
$ws.Cells.Item(16, 4).Value = 44160
$ws.Cells.Item(16, 10).Value = 62
$ws.Cells.Item(16, 11).Value = 19000
$ws.Cells.Item(16, 12).Value = 20000
$ws.Cells.Item(16, 13).Value = 19516
$ws.Cells.Item(16, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(16, 15).Value = 'Región del Maule'
$ws.Cells.Item(16, 16).Value = 781
$ws.Cells.Item(16, 17).Value = 25
# Row 17
$ws.Cells.Item(17, 4).Value = 44216
$ws.Cells.Item(17, 8).Value = 'Perfection'
$ws.Cells.Item(17, 10).Value = 43
$ws.Cells.Item(17, 11).Value = 24000
$ws.Cells.Item(17, 12).Value = 25000
$ws.Cells.Item(17, 13).Value = 24419
$ws.Cells.Item(17, 14).Value = '$/saco 30 kilos'
$ws.Cells.Item(17, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(17, 16).Value = 24419
$ws.Cells.Item(17, 17).Value = 1
# Row 18
$ws.Cells.Item(18, 4).Value = 44166
$ws.Cells.Item(18, 10).Value = 48
$ws.Cells.Item(18, 11).Value = 17000
$ws.Cells.Item(18, 12).Value = 18000
$ws.Cells.Item(18, 13).Value = 17479
$ws.Cells.Item(18, 15).Value = 'Región del Maule'
$ws.Cells.Item(18, 16).Value = 699
# Row 19
$ws.Cells.Item(19, 4).Value = 44215
$ws.Cells.Item(19, 10).Value = 42
$ws.Cells.Item(19, 13).Value = 19524
$ws.Cells.Item(19, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(19, 16).Value = 781
# Row 20
$ws.Cells.Item(20, 4).Value = 44186
$ws.Cells.Item(20, 10).Value = 30
$ws.Cells.Item(20, 11).Value = 21000
$ws.Cells.Item(20, 12).Value = 22000
$ws.Cells.Item(20, 13).Value = 21500
$ws.Cells.Item(20, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(20, 16).Value = 860
# Row 21
$ws.Cells.Item(21, 4).Value = 44273
$ws.Cells.Item(21, 10).Value = 22
$ws.Cells.Item(21, 11).Value = 20000
$ws.Cells.Item(21, 12).Value = 22000
$ws.Cells.Item(21, 13).Value = 21091
$ws.Cells.Item(21, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(21, 16).Value = 844
$ws.Cells.Item(21, 17).Value = 25
